# regen sval data to filter save games
# Update the per-appearance stat columns (TB, d2S, K, IP, sum) for rows 2-8.
# Column F (Win) is unchanged; column G (sum) = B+C+D+E and is stored as a
# literal value in the sheet, so it is updated alongside the inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2023-09-19
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.488907176552729

# Row 3 - 2023-08-19
$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.942654192641175

# Row 4 - 2023-07-16
$ws.Range("B4").Value = 0.06328177979961902
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.369310727790667

# Row 5 - 2023-07-09
$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.371470058157054

# Row 6 - 2023-07-06
$ws.Range("B6").Value = 0.3464964993005633
$ws.Range("C6").Value = 0.3375848360084654
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.336873824401267

# Row 7 - 2023-03-10
$ws.Range("B7").Value = 0.7287194209349384
$ws.Range("C7").Value = 0.05231270169004087
$ws.Range("D7").Value = 16.98373111632243
$ws.Range("E7").Value = 0.4998867070740569
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 18.26464994602146

# Row 8 - 2023-03-06
$ws.Range("B8").Value = 1.505614041169197
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 0.1529057820181812
$ws.Range("E8").Value = 0.4998867070740569
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 3.811642989160245
